# 1.1 Select Start.docx — apply the "enemy spawn zone fixed the menu
# funtionality" edit:
#   - Output Parameters: add "Provide Choice"
#   - Modules Called: drop the stray _GoBack bookmark
#   - Author: add "Delmis Spies"
#   - Date: add "12/2/2015"
#   - Peer Reviewer: add "Jan Cajas" (with the _GoBack bookmark moved
#     here, and a spell-check proofErr wrapper around "Cajas")

$d = $word.ActiveDocument

# Helper: replace a paragraph's content (everything up to, but not
# including, its end-of-paragraph mark) with an explicit run/bookmark/
# proofErr sequence supplied as literal WordprocessingML. Using
# Range.InsertXML wrapped in the minimal OOXML package envelope gives us
# real, independent <w:r> runs (and siblings like bookmarks/proofErr)
# instead of Word's usual "merge into the neighbouring run" behaviour.
function Set-ParaContent($para, [string]$innerXml) {
    $r = $d.Range($para.Range.Start, $para.Range.End)
    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>' +
           $innerXml +
           '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $r.InsertXML($pkg) | Out-Null
}

# Locate the table rows by their label text so this keeps working even
# if paragraph indices shift.
function Find-Paragraph([string]$startsWith) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.StartsWith($startsWith)) {
            return $p
        }
    }
    return $null
}

# NOTE: this interpreter's positional-parameter binding gets confused
# when a call's non-first argument is a parenthesized expression (even
# a literal like `('x')`), so every run-XML payload below is built into
# a plain variable first and that variable is passed to the helper.

# --- Output Parameters: Provide Choice --------------------------------
$pOutput = Find-Paragraph "Output Parameters:"
$xmlOutput = '<w:r><w:t xml:space="preserve">Output Parameters: </w:t></w:r>' +
    '<w:r><w:t>Provide Choice</w:t></w:r>'
Set-ParaContent $pOutput $xmlOutput

# --- Modules Called: strip the leftover _GoBack bookmark --------------
$pModules = Find-Paragraph "Modules Called:"
$xmlModules = '<w:r><w:t xml:space="preserve">Modules Called: </w:t></w:r>' +
    '<w:r><w:t>Select New Game 1.1.1, Select Load Game 1.1.2</w:t></w:r>'
Set-ParaContent $pModules $xmlModules

# --- Author: Delmis Spies ----------------------------------------------
$pAuthor = Find-Paragraph "Author:"
$xmlAuthor = '<w:r><w:t xml:space="preserve">Author: </w:t></w:r>' +
    '<w:r><w:t>Delmis Spies</w:t></w:r>'
Set-ParaContent $pAuthor $xmlAuthor

# --- Date: 12/2/2015 -----------------------------------------------------
$pDate = Find-Paragraph "Date:"
$xmlDate = '<w:r><w:t xml:space="preserve">Date:  </w:t></w:r>' +
    '<w:r><w:t>12/2/2015</w:t></w:r>'
Set-ParaContent $pDate $xmlDate

# --- Peer Reviewer: Jan Cajas (+ _GoBack bookmark, + spell-check proof) -
$pPeer = Find-Paragraph "Peer Reviewer"
$xmlPeer = '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:r><w:t>Peer Reviewer</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">Jan </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Cajas</w:t></w:r>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:proofErr w:type="spellEnd"/>'
Set-ParaContent $pPeer $xmlPeer
